# Updates the weekly Sandia (watermelon) price records (rows 3-14) of
# "Mapocho Venta Directa de Santiago" with the new reported values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 44194
$ws.Range("I3").Value = "Extra"
$ws.Range("J3").Value = 120
$ws.Range("K3").Value = 3500
$ws.Range("L3").Value = 3500
$ws.Range("M3").Value = 3500
$ws.Range("O3").Value = "Región de O'Higgins"
$ws.Range("P3").Value = 3500
$ws.Range("D4").Value = 44194
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 3000
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = 3000
$ws.Range("P4").Value = 3000
$ws.Range("D5").Value = 44223
$ws.Range("H5").Value = "Americana O Klondike"
$ws.Range("I5").Value = "Extra"
$ws.Range("J5").Value = 340
$ws.Range("K5").Value = 2500
$ws.Range("L5").Value = 2500
$ws.Range("M5").Value = 2500
$ws.Range("P5").Value = 2500
$ws.Range("D6").Value = 44223
$ws.Range("H6").Value = "Americana O Klondike"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 2000
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = 2000
$ws.Range("P6").Value = 2000
$ws.Range("D7").Value = 44223
$ws.Range("H7").Value = "Americana O Klondike"
$ws.Range("I7").Value = "Segunda"
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 1500
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = 1500
$ws.Range("P7").Value = 1500
$ws.Range("I8").Value = "Tercera"
$ws.Range("J8").Value = 160
$ws.Range("K8").Value = 1000
$ws.Range("L8").Value = 1000
$ws.Range("M8").Value = 1000
$ws.Range("P8").Value = 1000
$ws.Range("D9").Value = 44167
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("K9").Value = 5000
$ws.Range("L9").Value = 5000
$ws.Range("M9").Value = 5000
$ws.Range("P9").Value = 5000
$ws.Range("D10").Value = 44167
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("J10").Value = 560
$ws.Range("K10").Value = 3000
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = 3000
$ws.Range("P10").Value = 3000
$ws.Range("D11").Value = 44167
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("J11").Value = 450
$ws.Range("K11").Value = 2000
$ws.Range("L11").Value = 2000
$ws.Range("M11").Value = 2000
$ws.Range("P11").Value = 2000
$ws.Range("D12").Value = 44217
$ws.Range("I12").Value = "Extra"
$ws.Range("K12").Value = 2500
$ws.Range("L12").Value = 2500
$ws.Range("M12").Value = 2500
$ws.Range("P12").Value = 2500
$ws.Range("D13").Value = 44217
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 280
$ws.Range("K13").Value = 2000
$ws.Range("L13").Value = 2000
$ws.Range("M13").Value = 2000
$ws.Range("P13").Value = 2000
$ws.Range("D14").Value = 44312
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 180
$ws.Range("K14").Value = 2500
$ws.Range("L14").Value = 2500
$ws.Range("M14").Value = 2500
$ws.Range("O14").Value = "Perú"
$ws.Range("P14").Value = 2500
